# Update the "想去人数" (column F) figures on all four sheets to the
# values published at commit 456a3b4 (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 3568
$ws.Cells.Item(5, 6).Value = 8332
$ws.Cells.Item(7, 6).Value = 127
$ws.Cells.Item(8, 6).Value = 2230
$ws.Cells.Item(10, 6).Value = 100
$ws.Cells.Item(11, 6).Value = 75
$ws.Cells.Item(12, 6).Value = 653
$ws.Cells.Item(13, 6).Value = 108
$ws.Cells.Item(14, 6).Value = 7471
$ws.Cells.Item(15, 6).Value = 450
$ws.Cells.Item(16, 6).Value = 7711
$ws.Cells.Item(18, 6).Value = 57805
$ws.Cells.Item(19, 6).Value = 57805
$ws.Cells.Item(20, 6).Value = 4809
$ws.Cells.Item(21, 6).Value = 1061
$ws.Cells.Item(22, 6).Value = 953
$ws.Cells.Item(23, 6).Value = 504
$ws.Cells.Item(24, 6).Value = 112
$ws.Cells.Item(25, 6).Value = 930
$ws.Cells.Item(27, 6).Value = 618
$ws.Cells.Item(28, 6).Value = 5311
$ws.Cells.Item(29, 6).Value = 602
$ws.Cells.Item(30, 6).Value = 119
$ws.Cells.Item(31, 6).Value = 51
$ws.Cells.Item(32, 6).Value = 915
$ws.Cells.Item(33, 6).Value = 1401
$ws.Cells.Item(34, 6).Value = 1979
$ws.Cells.Item(35, 6).Value = 21
$ws.Cells.Item(36, 6).Value = 187
$ws.Cells.Item(37, 6).Value = 233
$ws.Cells.Item(38, 6).Value = 1092
$ws.Cells.Item(40, 6).Value = 731
$ws.Cells.Item(42, 6).Value = 785
$ws.Cells.Item(43, 6).Value = 271
$ws.Cells.Item(44, 6).Value = 220
$ws.Cells.Item(47, 6).Value = 206
$ws.Cells.Item(50, 6).Value = 2489

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 203
$ws.Cells.Item(3, 6).Value = 37
$ws.Cells.Item(6, 6).Value = 144
$ws.Cells.Item(9, 6).Value = 7646
$ws.Cells.Item(14, 6).Value = 5
$ws.Cells.Item(22, 6).Value = 34
$ws.Cells.Item(23, 6).Value = 40
$ws.Cells.Item(26, 6).Value = 131
$ws.Cells.Item(29, 6).Value = 3
$ws.Cells.Item(32, 6).Value = 3
$ws.Cells.Item(45, 6).Value = 31
$ws.Cells.Item(48, 6).Value = 280

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 2376
$ws.Cells.Item(5, 6).Value = 1613
$ws.Cells.Item(7, 6).Value = 683
$ws.Cells.Item(8, 6).Value = 2424
$ws.Cells.Item(9, 6).Value = 9442
$ws.Cells.Item(10, 6).Value = 1771
$ws.Cells.Item(15, 6).Value = 273
$ws.Cells.Item(16, 6).Value = 2395
$ws.Cells.Item(17, 6).Value = 49
$ws.Cells.Item(19, 6).Value = 522

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 3568
$ws.Cells.Item(3, 6).Value = 2376
$ws.Cells.Item(5, 6).Value = 8332
$ws.Cells.Item(6, 6).Value = 683
$ws.Cells.Item(8, 6).Value = 127
$ws.Cells.Item(9, 6).Value = 273
$ws.Cells.Item(10, 6).Value = 75
$ws.Cells.Item(11, 6).Value = 653
$ws.Cells.Item(12, 6).Value = 108
$ws.Cells.Item(13, 6).Value = 7472
$ws.Cells.Item(14, 6).Value = 7711
$ws.Cells.Item(15, 6).Value = 57805
$ws.Cells.Item(16, 6).Value = 203
$ws.Cells.Item(17, 6).Value = 37
$ws.Cells.Item(18, 6).Value = 4809
$ws.Cells.Item(19, 6).Value = 1061
$ws.Cells.Item(20, 6).Value = 953
$ws.Cells.Item(21, 6).Value = 504
$ws.Cells.Item(22, 6).Value = 618
$ws.Cells.Item(23, 6).Value = 144
$ws.Cells.Item(24, 6).Value = 5311
$ws.Cells.Item(25, 6).Value = 602
$ws.Cells.Item(26, 6).Value = 119
$ws.Cells.Item(27, 6).Value = 915
$ws.Cells.Item(28, 6).Value = 1401
$ws.Cells.Item(30, 6).Value = 522
$ws.Cells.Item(31, 6).Value = 5
$ws.Cells.Item(33, 6).Value = 187
$ws.Cells.Item(35, 6).Value = 731
$ws.Cells.Item(36, 6).Value = 785
$ws.Cells.Item(37, 6).Value = 271
$ws.Cells.Item(39, 6).Value = 40
$ws.Cells.Item(43, 6).Value = 3
$ws.Cells.Item(45, 6).Value = 206
$ws.Cells.Item(48, 6).Value = 44
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(50, 6).Value = 31
